$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-05 Sunday" "2025-10-06 Monday"

Replace-Text "701×5=3505" "872×6=5232"
Replace-Text "970×2=1940" "496×6=2976"
Replace-Text "151×9=1359" "811×5=4055"
Replace-Text "991×7=6937" "368×6=2208"
Replace-Text "172×4=688" "159×2=318"
Replace-Text "943×2=1886" "236×4=944"
Replace-Text "592×6=3552" "421×5=2105"
Replace-Text "398×5=1990" "252×6=1512"
Replace-Text "637×7=4459" "471×7=3297"
Replace-Text "665×2=1330" "148×6=888"
Replace-Text "817×7=5719" "234×2=468"
Replace-Text "950×5=4750" "718×2=1436"
Replace-Text "908×4=3632" "848×2=1696"
Replace-Text "731×9=6579" "163×3=489"
Replace-Text "432×6=2592" "894×6=5364"
Replace-Text "999×6=5994" "971×9=8739"
Replace-Text "616×7=4312" "591×7=4137"
Replace-Text "936×8=7488" "502×9=4518"
Replace-Text "870×8=6960" "780×6=4680"
Replace-Text "232×8=1856" "209×8=1672"
Replace-Text "423×4=1692" "579×7=4053"
Replace-Text "766×9=6894" "175×6=1050"
Replace-Text "231×6=1386" "499×7=3493"
Replace-Text "937×9=8433" "492×2=984"
Replace-Text "489×7=3423" "800×8=6400"
